$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "DANH SÁCH NỢ" (sheet1) -- add two new debt entries (rows 15 & 16),
# which were pre-existing blank rows, and insert four more blank rows after
# row 17 so the table grows from 16 to 20 data rows (STT 1..20).
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("DANH SÁCH NỢ")

# Insert 4 blank rows right after row 17 (pushes old row 18.. down to 22..)
$wsDebt.Rows("18:21").Insert()

# The newly inserted rows 18:21 land with generic default formatting; copy
# the normal data-row look from row 16 into rows 18-20 ...
$wsDebt.Range("A16:M16").Copy()
$wsDebt.Range("A18:M20").PasteSpecial(-4122)

# ...and move the special "closing border" formatting that used to belong to
# row 17 down onto the new last row (21) ...
$wsDebt.Range("A17:M17").Copy()
$wsDebt.Range("A21:M21").PasteSpecial(-4122)

# ...while row 17 itself becomes just another normal data row now.
$wsDebt.Range("A16:M16").Copy()
$wsDebt.Range("A17:M17").PasteSpecial(-4122)
$wsDebt.Application.CutCopyMode = $false

# Re-number the STT column for the (now blank) trailing rows.
$wsDebt.Range("A18").Value = 17
$wsDebt.Range("A19").Value = 18
$wsDebt.Range("A20").Value = 19
$wsDebt.Range("A21").Value = 20

# Fill in the two new debt records (rows 15 & 16).
$wsDebt.Range("B15").Value = "Huỳnh Quốc Phú"
$wsDebt.Range("C15").Value = "Nạp quân huy"
$wsDebt.Range("D15").Value = 75000
$wsDebt.Range("E15").Value = 0
$wsDebt.Range("F15").Formula = "=(D15+I15)-E15"
$wsDebt.Range("G15").Value = 0
$wsDebt.Range("H15").Value = 0
$wsDebt.Range("I15").Formula = "=D15*H15"
$wsDebt.Range("J15").Value = 46019
$wsDebt.Range("K15").Value = 46025
$wsDebt.Range("M15").Value = "Chưa trả đủ"

$wsDebt.Range("B16").Value = "Bùi Anh Tài"
$wsDebt.Range("C16").Value = "Nạp Robux"
$wsDebt.Range("D16").Value = 100000
$wsDebt.Range("E16").Value = 0
$wsDebt.Range("F16").Formula = "=(D16+I16)-E16"
$wsDebt.Range("G16").Value = 0
$wsDebt.Range("H16").Value = 0
$wsDebt.Range("I16").Formula = "=D16*H16"
$wsDebt.Range("J16").Value = 46019
$wsDebt.Range("K16").Value = 46025
$wsDebt.Range("M16").Value = "Chưa trả đủ"

# Extend the three totals formulas (now 4 rows further down) to cover the
# full, longer data range.
$wsDebt.Range("F23").Formula = "=SUM(D2:D21)"
$wsDebt.Range("F24").Formula = "=SUM(E2:E21)"
$wsDebt.Range("F25").Formula = "=-SUM(F2:F21)"

# Move the "Điều khoản dịch vụ:" hyperlink from its old cell (D18) onto its
# new cell (D22).
foreach ($h in $wsDebt.Hyperlinks) {
    $h.Delete()
}
$wsDebt.Hyperlinks.Add($wsDebt.Range("D22"), "https://tinyurl.com/dieukhoan29")

# Grow the AutoFilter range to match the new table extent.
$wsDebt.AutoFilterMode = $false
$wsDebt.Range("A1:M25").AutoFilter()

# Keep the hidden _FilterDatabase defined name for this sheet in sync too.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*" -and $n.RefersTo -like "*DANH S*") {
        $n.RefersTo = "='DANH SÁCH NỢ'!`$A`$1:`$M`$25"
    }
}

# ---------------------------------------------------------------------------
# Sheet "THONG KE NAP " (sheet2) -- log the same two payments in the running
# deposit/expense ledger (rows 116 & 117 were blank placeholders).
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("THONG KE NAP ")

$wsLog.Range("A116").Value = 46019
$wsLog.Range("B116").Value = "Huỳnh Quốc Phú"
$wsLog.Range("C116").Value = 75000
$wsLog.Range("D116").Value = "Nạp quân huy"

$wsLog.Range("A117").Value = 46019
$wsLog.Range("B117").Value = "Bùi Anh Tài"
$wsLog.Range("C117").Value = 100000
$wsLog.Range("D117").Value = "Nạp Robux"

# ---------------------------------------------------------------------------
# Restore the on-screen selections shown in the saved workbook. Do the
# non-active sheet first so the originally active sheet (DANH SÁCH NỢ) ends
# up selected last / still the active tab.
# ---------------------------------------------------------------------------
$wsLog.Range("D118").Select()
$wsDebt.Range("O12").Select()
